# Weekly update: insert two new daily records (row date 2021-09-20 / serial 44459)
# for "Jengibre" at Vega Central Mapocho de Santiago, right after the existing
# row 19 (date 44431). This pushes the previously existing rows 20-55 down to
# 22-57, which naturally reproduces the rest of the diff (every later row's
# values simply shift down by two rows, with no value changes of their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before current row 20 (shifts old rows 20:55 -> 22:57)
$ws.Rows("20:21").Insert()

# --- New row 20: Primera ---
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C20").Value = 'Metropolitana'
$ws.Range("D20").Value = 44459
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100114007
$ws.Range("G20").Value = 'Jengibre'
$ws.Range("H20").Value = 'Sin especificar'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 970
$ws.Range("K20").Value = 13000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13495
$ws.Range("N20").Value = '$/caja 13 kilos'
$ws.Range("O20").Value = 'Perú'
$ws.Range("P20").Value = 1038
$ws.Range("Q20").Value = 13
$ws.Range("R20").Value = 'Hortaliza'

# --- New row 21: Segunda ---
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C21").Value = 'Metropolitana'
$ws.Range("D21").Value = 44459
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100114007
$ws.Range("G21").Value = 'Jengibre'
$ws.Range("H21").Value = 'Sin especificar'
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 520
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 11500
$ws.Range("N21").Value = '$/caja 13 kilos'
$ws.Range("O21").Value = 'Perú'
$ws.Range("P21").Value = 885
$ws.Range("Q21").Value = 13
$ws.Range("R21").Value = 'Hortaliza'
